# "Added Update profile scripts"
# Adds a new "Expected Title" column (E) to the existing credential test
# matrix and a new block of rows describing an "Update profile" test case
# (short values in row 6, first/last name in rows 7-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

# New column E: expected outcome for each existing test row.
$ws.Range("E1").Value = "Expeceted Title"
$ws.Range("E2").Value = "Register"

# New "Update profile" test data block (row 5 left blank as a separator).
$ws.Range("A6").Value = "r"
$ws.Range("B6").Value = "Rj"
$ws.Range("A7").Value = "Rutu"
$ws.Range("A8").Value = "Jadhav"

$ws.Range("A8").Select()
